$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a flat daily price log for "Cilantro" at Femacal de La Calera.
# A new daily record (2023-09-04) is inserted before the current row 522,
# pushing the existing row 522..634 records down to 523..635.
$ws.Rows.Item(522).Insert()

$ws.Range("A522").Value = 3
$ws.Range("B522").Value = "Femacal de La Calera"
$ws.Range("C522").Value = "Coquimbo"
$ws.Range("D522").Value = "2023-09-04"
$ws.Range("E522").Value = 5
$ws.Range("F522").Value = 100112040
$ws.Range("G522").Value = "Cilantro"
$ws.Range("H522").Value = "Sin especificar"
$ws.Range("I522").Value = "Primera"
$ws.Range("J522").Value = 170
$ws.Range("K522").Value = 4000
$ws.Range("L522").Value = 4300
$ws.Range("M522").Value = 4106
$ws.Range("N522").Value = "$/docena de atados (3 kilos)"
$ws.Range("O522").Value = "Provincia de Quillota"
$ws.Range("P522").Value = 1369
$ws.Range("Q522").Value = 3
$ws.Range("R522").Value = "Hortaliza"
